$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "69.641.67"
$ws.Range("D3").Formula = "3.689.42"
$ws.Range("E3").Formula = "  +0.43%  "
$ws.Range("E4").Formula = "  -0.07%  "
$ws.Range("D5").Formula = "'667.38"
$ws.Range("E5").Formula = "  -1.01%  "
$ws.Range("D6").Formula = "'160.05"
$ws.Range("E6").Formula = "  +1.57%  "
$ws.Range("D7").Formula = "'1.00"
$ws.Range("E7").Formula = "  +0.00%  "
$ws.Range("E9").Formula = "  +0.37%  "
$ws.Range("D10").Formula = "'7.12"
$ws.Range("E10").Formula = "  +3.84%  "
$ws.Range("E11").Formula = "  +1.79%  "
$ws.Range("D12").Formula = "'0.0000233"
$ws.Range("E12").Formula = "  +1.26%  "
$ws.Range("D13").Formula = "'32.80"
$ws.Range("E13").Formula = "  +1.91%  "
$ws.Range("D14").Formula = "3.718.10"
$ws.Range("E14").Formula = "  +1.13%  "
$ws.Range("D15").Formula = "69.658.87"
$ws.Range("E15").Formula = "  +0.65%  "
$ws.Range("D16").Formula = "'0.118"
$ws.Range("E16").Formula = "  +2.62%  "
$ws.Range("D17").Formula = "'16.14"
$ws.Range("E17").Formula = "  +1.25%  "
$ws.Range("D18").Formula = "'6.46"
$ws.Range("E18").Formula = "  +1.04%  "
$ws.Range("D19").Formula = "'469.76"
$ws.Range("E19").Formula = "  +0.59%  "
$ws.Range("E20").Formula = "  -1.36%  "
$ws.Range("D21").Formula = "'0.645"
$ws.Range("E21").Formula = "  -0.17%  "
$ws.Range("D22").Formula = "'79.89"
$ws.Range("E22").Formula = "  +0.26%  "
$ws.Range("D23").Formula = "3.836.56"
$ws.Range("E23").Formula = "  +0.45%  "
$ws.Range("E24").Formula = "  +5.69%  "
$ws.Range("E25").Formula = "  -0.01%  "
$ws.Range("D26").Formula = "'10.90"
$ws.Range("E26").Formula = "  +0.73%  "
$ws.Range("D27").Formula = "'9.04"
$ws.Range("E27").Formula = "  +0.57%  "
$ws.Range("D28").Formula = "'2.68"
$ws.Range("E28").Formula = "  +0.41%  "
$ws.Range("E29").Formula = "  -2.39%  "
# Rows 30/31 swap: ImmutableX <-> Binance-PegBSC-USD
$ws.Range("B30").Formula = "Binance-PegBSC-USD"
$ws.Range("C30").Formula = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Formula = "'1.01"
$ws.Range("E30").Formula = "  +0.78%  "
$ws.Range("B31").Formula = "ImmutableX"
$ws.Range("C31").Formula = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Formula = "'2.00"
$ws.Range("E31").Formula = "  +1.71%  "
$ws.Range("D32").Formula = "'26.75"
$ws.Range("E32").Formula = "  -0.40%  "
$ws.Range("E33").Formula = "  +3.64%  "
$ws.Range("D34").Formula = "'6.47"
$ws.Range("E34").Formula = "  -1.58%  "
$ws.Range("D35").Formula = "3.683.01"
$ws.Range("E35").Formula = "  +0.38%  "
$ws.Range("E36").Formula = "  +4.05%  "
$ws.Range("D37").Formula = "'6.09"
$ws.Range("E37").Formula = "  -1.54%  "
$ws.Range("E38").Formula = "  -0.03%  "
$ws.Range("E39").Formula = "  +1.64%  "
# Rows 40/41 swap: FirstDigitalUSD <-> Monero
$ws.Range("B40").Formula = "Monero"
$ws.Range("C40").Formula = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Formula = "'178.12"
$ws.Range("E40").Formula = "  +3.14%  "
$ws.Range("B41").Formula = "FirstDigitalUSD"
$ws.Range("C41").Formula = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Formula = "'1.00"
$ws.Range("E41").Formula = "  -0.04%  "
$ws.Range("D42").Formula = "'0.0904"
$ws.Range("E42").Formula = "  +1.02%  "
$ws.Range("D43").Formula = "'0.934"
$ws.Range("E43").Formula = "  -0.61%  "
$ws.Range("D44").Formula = "'47.02"
$ws.Range("E44").Formula = "  -1.31%  "
$ws.Range("E45").Formula = "  +4.03%  "
$ws.Range("D46").Formula = "'27.54"
$ws.Range("E46").Formula = "  -1.63%  "
$ws.Range("E47").Formula = "  +0.38%  "
$ws.Range("D48").Formula = "'0.000272"
$ws.Range("E48").Formula = "  -0.78%  "
$ws.Range("E49").Formula = "  -0.26%  "
$ws.Range("D50").Formula = "'7.84"
$ws.Range("E50").Formula = "  +1.22%  "
$ws.Range("E51").Formula = "  +0.28%  "
